$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted above the current row 349,
# pushing all subsequent rows (349:425) down by one (to 350:426).
$ws.Rows("349:349").Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A349").Value = 6
$ws.Range("B349").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C349").Value = "Metropolitana"
$ws.Range("D349").Value = 45135
$ws.Range("E349").Value = 13
$ws.Range("F349").Value = 100112026
$ws.Range("G349").Value = "Haba"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 330
$ws.Range("K349").Value = 13000
$ws.Range("L349").Value = 14000
$ws.Range("M349").Value = 13636
$ws.Range("N349").Value = "$/saco 25 kilos"
$ws.Range("O349").Value = "Región de Coquimbo"
$ws.Range("P349").Value = 545
$ws.Range("Q349").Value = 25
$ws.Range("R349").Value = "Hortaliza"
